# Regenerate the K (strikeout) column (column G) for the save_data sheet.
# This mirrors the data pipeline's re-scrape/regen step that replaced the
# old "Strike#" derived values with the actual K (strikeouts) values,
# after recalculating std/mean and writing the resulting s_vals back
# into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2-52 (one per game row), in top-to-bottom order.
$kValues = @(2,1,1,0,0,0,1,0,0,0,2,1,2,0,0,1,1,0,2,2,2,0,1,2,0,2,0,0,2,1,0,0,2,0,2,3,0,1,1,1,0,1,1,1,2,0,2,0,2,1,1)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
